$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number + report week date range) ---
$ws.Range("A8").Value = "Volume 31   Number  27"
$ws.Range("C9").Value = "Report Covering the Week  7/1/2024  Through  7/7/2024"

# --- C15 / C27 / C28 switch from a numeric count to the text placeholder "0" ---
# (mirrors the existing "0"-placeholder cells, e.g. D15/D27/D28, so the style/shared-string
#  entry is reused exactly instead of minting a new one)
$ws.Range("D15").Copy($ws.Range("C15"))
$ws.Range("D27").Copy($ws.Range("C27"))
$ws.Range("D28").Copy($ws.Range("C28"))

# --- Updated weekly crime statistics ---
$ws.Range("N15").Value = -42.857142857142
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 0
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = -50
$ws.Range("I16").Value = 30
$ws.Range("J16").Value = 35
$ws.Range("K16").Value = -14.285714285714
$ws.Range("L16").Value = -30.232558139534
$ws.Range("M16").Value = -31.818181818181
$ws.Range("N16").Value = -82.142857142857
$ws.Range("D17").Value = 3
$ws.Range("F17").Value = 6
$ws.Range("G17").Value = 7
$ws.Range("H17").Value = -14.285714285714
$ws.Range("J17").Value = 53
$ws.Range("K17").Value = -5.660377358490
$ws.Range("L17").Value = 25
$ws.Range("M17").Value = 78.571428571428
$ws.Range("N17").Value = -23.076923076923
$ws.Range("C18").Value = 7
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = 16.666666666666
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = 13.333333333333
$ws.Range("I18").Value = 133
$ws.Range("J18").Value = 165
$ws.Range("K18").Value = -19.393939393939
$ws.Range("L18").Value = -5
$ws.Range("M18").Value = 8.130081300813
$ws.Range("N18").Value = -75.046904315197
$ws.Range("C19").Value = 4
$ws.Range("E19").Value = -66.666666666666
$ws.Range("F19").Value = 30
$ws.Range("G19").Value = 45
$ws.Range("H19").Value = -33.333333333333
$ws.Range("I19").Value = 237
$ws.Range("J19").Value = 341
$ws.Range("K19").Value = -30.498533724340
$ws.Range("L19").Value = -19.112627986348
$ws.Range("M19").Value = 32.402234636871
$ws.Range("N19").Value = -10.902255639097
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 400
$ws.Range("F20").Value = 15
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = 50
$ws.Range("I20").Value = 130
$ws.Range("J20").Value = 75
$ws.Range("K20").Value = 73.333333333333
$ws.Range("L20").Value = 165.30612244898
$ws.Range("M20").Value = 78.082191780821
$ws.Range("N20").Value = -92.243436754176
$ws.Range("C21").Value = 17
$ws.Range("D21").Value = 23
$ws.Range("E21").Value = -26.086956521739
$ws.Range("F21").Value = 73
$ws.Range("G21").Value = 87
$ws.Range("H21").Value = -16.091954022988
$ws.Range("I21").Value = 584
$ws.Range("J21").Value = 677
$ws.Range("K21").Value = -13.737075332348
$ws.Range("L21").Value = 2.998236331569
$ws.Range("M21").Value = 29.490022172949
$ws.Range("N21").Value = -78.497790868924
$ws.Range("C24").Value = 8
$ws.Range("D24").Value = 16
$ws.Range("E24").Value = -50
$ws.Range("F24").Value = 38
$ws.Range("G24").Value = 60
$ws.Range("H24").Value = -36.666666666666
$ws.Range("I24").Value = 286
$ws.Range("J24").Value = 295
$ws.Range("K24").Value = -3.050847457627
$ws.Range("L24").Value = -29.207920792079
$ws.Range("M24").Value = 16.734693877551
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = -33.333333333333
$ws.Range("F25").Value = 11
$ws.Range("G25").Value = 18
$ws.Range("H25").Value = -38.888888888888
$ws.Range("I25").Value = 60
$ws.Range("J25").Value = 67
$ws.Range("K25").Value = -10.447761194029
$ws.Range("L25").Value = 5.263157894736
$ws.Range("C26").Value = 7
$ws.Range("E26").Value = 250
$ws.Range("F26").Value = 25
$ws.Range("G26").Value = 13
$ws.Range("H26").Value = 92.307692307692
$ws.Range("I26").Value = 111
$ws.Range("J26").Value = 112
$ws.Range("K26").Value = -0.892857142857
$ws.Range("L26").Value = -12.598425196850
$ws.Range("M26").Value = 32.142857142857
$ws.Range("L28").Value = -57.142857142857
$ws.Range("L31").Value = 0
